$d = $word.ActiveDocument

# 1 & 4) Strip the stray trailing newline that follows "casebook." at the end
#        of the two plain narrative paragraphs ("...annotatable resource in
#        the casebook." and "...second chapter of the casebook."). Both
#        paragraphs consist of a single run, so a direct, scoped text
#        replacement only touches that run.
$d.Content.Find.Execute("casebook.`n", $true, $false, $false, $false, $false, `
    $true, 1, $false, "casebook.", 2) | Out-Null

# 2) "...elided: [ … ];\nreplaced: foo bar baz..." -> collapse the embedded
#    newline between ";" and "replaced: " into a single space, touching only
#    that one run (bounded range, so neighboring runs such as the Elision
#    styled "[ … ]" run are left alone).
$search1 = $d.Content
$found1 = $search1.Find.Execute(";`nreplaced: ")
if ($found1) {
    $nl1 = $search1.Start + 1
    $rng1 = $d.Range($nl1, $nl1 + 1)
    $rng1.Text = " "
}

# 3) "...content to link; noted:\ncontent to note;..." -> collapse the
#    embedded newline after "noted:" into a single space, again touching
#    only the single character so the neighboring hyperlink run keeps its
#    own formatting.
$search2 = $d.Content
$found2 = $search2.Find.Execute("; noted:`n")
if ($found2) {
    $nl2 = $search2.End - 1
    $rng2 = $d.Range($nl2, $nl2 + 1)
    $rng2.Text = " "
}
